$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.220.88"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.286.40"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "2.285.98"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0949"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("E12").Value = "  +4.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "2.689.91"
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").Value = "54.206.37"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "2.287.23"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "301.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").Value = "2.379.82"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.149"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.74%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "0.0₃0688"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.874"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.374"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.544"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "238.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0484"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.47%  "
